$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "[캐글] 성인 인구조사 소득 예측 대회 커널"
$ws.Range("E4").Value = "https://teddylee777.github.io/kaggle/kakr4"

$ws.Range("D12").Value = "TensorFlow 2.7.0이 릴리스되었습니다."
$ws.Range("E12").Value = "https://tensorflow.blog/2021/11/06/tensorflow-2-7-0%ec%9d%b4-%eb%a6%b4%eb%a6%ac%ec%8a%a4%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D20").Value = "[머신러닝 기초] 지도학습 - classification (Logistic Regression)"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/588"

$ws.Range("D44").Value = "Non-Fungible Token (NFT)의 개념과 견해"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/105"
